$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = "SELECT`n    smp.sample_id AS ""Sample ID"",`n    prt.participant_id AS ""Participant ID"",`n    std.dbgap_accession AS ""Study ID"",`n    smp.anatomic_site AS ""Sample Anatomic Site"",`n    COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS ""Age at Sample Collection (days)"",`n    COALESCE(smp.sample_tumor_status, '') AS ""Sample Tumor Status"",`n    COALESCE(smp.tumor_classification, '') AS ""Sample Tumor Classification"",`n    Null  AS ""Sample Diagnosis""`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.""study.id""`nLEFT JOIN `n    df_sample smp ON prt.id = smp.""participant.id""`nLEFT JOIN `n    df_diagnosis dgn ON prt.id = dgn.""participant.id""`nLEFT JOIN `n    df_survival srv ON prt.id = srv.""participant.id""`nWHERE `n    std.dbgap_accession = 'phs002371' `n    AND smp.anatomic_site = 'C42.1 : Bone marrow'`n`tand srv.last_known_survival_status ='Alive'`nORDER BY `n    smp.sample_id ASC;"

$ws.Range("B4").Value = $newQuery
$ws.Rows.Item(4).RowHeight = 409.5

$ws.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
